$d = $word.ActiveDocument

$anchorText = "Expect your requirements will change frequently as you refine your understanding of the customer needs and your technology stack."
$newText = "Write as many requirements as you actually need, you don’t need to write every single requirement in one shot."

# Locate the anchor paragraph's index within the document.
$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $anchorText) {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not locate the anchor paragraph."
}

$anchor = $d.Paragraphs.Item($anchorIndex)
$anchor.Range.InsertParagraphAfter()

# Re-fetch the freshly created (now-live) paragraph and fill in its text,
# preserving the style/numbering it inherited from the anchor paragraph.
$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newPara.Range.Text = $newText
